# Update NATMI LR-pairs output (Tnfsf13b -> Tnfrsf13c) with refreshed TPM
# figures. The sending/target cluster pairing for the FAPs / MuSCs /
# Resolving-Mac senders changes, the "ECs" sending-cluster row disappears,
# and two additional target rows (ECs as target for FAPs/MuSCs/Resolving-Mac
# senders) are appended, growing the table from 4 to 6 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Tnfsf13b/Tnfrsf13c -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Tnfsf13b"
$ws.Range("C2").Value = "Tnfrsf13c"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.643401666666666
$ws.Range("H2").Value = 4.930204999999999
$ws.Range("I2").Value = 0.6433537405875911
$ws.Range("J2").Value = 0.6433537405875911
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3064233333333333
$ws.Range("N2").Value = 0.91927
$ws.Range("O2").Value = 0.303622416540745
$ws.Range("P2").Value = 0.303622416540745
$ws.Range("Q2").Value = 0.5035766167055554
$ws.Range("R2").Value = 4.532189550349999
$ws.Range("S2").Value = 0.1953366174077319
$ws.Range("T2").Value = 0.1953366174077319

# Row 3: FAPs -> Tnfsf13b/Tnfrsf13c -> MuSCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Tnfsf13b"
$ws.Range("C3").Value = "Tnfrsf13c"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.643401666666666
$ws.Range("H3").Value = 4.930204999999999
$ws.Range("I3").Value = 0.6433537405875911
$ws.Range("J3").Value = 0.6433537405875911
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7028016666666667
$ws.Range("N3").Value = 2.108405
$ws.Range("O3").Value = 0.696377583459255
$ws.Range("P3").Value = 0.696377583459255
$ws.Range("Q3").Value = 1.154985430336111
$ws.Range("R3").Value = 10.394868873025
$ws.Range("S3").Value = 0.4480171231798591
$ws.Range("T3").Value = 0.4480171231798591

# Row 4: MuSCs -> Tnfsf13b/Tnfrsf13c -> ECs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Tnfsf13b"
$ws.Range("C4").Value = "Tnfrsf13c"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.08257399999999999
$ws.Range("H4").Value = 0.247722
$ws.Range("I4").Value = 0.03232581106177922
$ws.Range("J4").Value = 0.03232581106177923
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3064233333333333
$ws.Range("N4").Value = 0.91927
$ws.Range("O4").Value = 0.303622416540745
$ws.Range("P4").Value = 0.303622416540745
$ws.Range("Q4").Value = 0.02530260032666666
$ws.Range("R4").Value = 0.22772340294
$ws.Range("S4").Value = 0.009814840871216953
$ws.Range("T4").Value = 0.009814840871216955

# Row 5: MuSCs -> Tnfsf13b/Tnfrsf13c -> MuSCs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Tnfsf13b"
$ws.Range("C5").Value = "Tnfrsf13c"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.08257399999999999
$ws.Range("H5").Value = 0.247722
$ws.Range("I5").Value = 0.03232581106177922
$ws.Range("J5").Value = 0.03232581106177923
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7028016666666667
$ws.Range("N5").Value = 2.108405
$ws.Range("O5").Value = 0.696377583459255
$ws.Range("P5").Value = 0.696377583459255
$ws.Range("Q5").Value = 0.05803314482333333
$ws.Range("R5").Value = 0.5222983034099999
$ws.Range("S5").Value = 0.02251097019056227
$ws.Range("T5").Value = 0.02251097019056228

# Row 6 (new): Resolving-Mac -> Tnfsf13b/Tnfrsf13c -> ECs
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Tnfsf13b"
$ws.Range("C6").Value = "Tnfrsf13c"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.8284536666666668
$ws.Range("H6").Value = 2.485361
$ws.Range("I6").Value = 0.3243204483506297
$ws.Range("J6").Value = 0.3243204483506297
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3064233333333333
$ws.Range("N6").Value = 0.91927
$ws.Range("O6").Value = 0.303622416540745
$ws.Range("P6").Value = 0.303622416540745
$ws.Range("Q6").Value = 0.2538575340522222
$ws.Range("R6").Value = 2.28471780647
$ws.Range("S6").Value = 0.09847095826179604
$ws.Range("T6").Value = 0.09847095826179605

# Row 7 (new): Resolving-Mac -> Tnfsf13b/Tnfrsf13c -> MuSCs
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Tnfsf13b"
$ws.Range("C7").Value = "Tnfrsf13c"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.8284536666666668
$ws.Range("H7").Value = 2.485361
$ws.Range("I7").Value = 0.3243204483506297
$ws.Range("J7").Value = 0.3243204483506297
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7028016666666667
$ws.Range("N7").Value = 2.108405
$ws.Range("O7").Value = 0.696377583459255
$ws.Range("P7").Value = 0.696377583459255
$ws.Range("Q7").Value = 0.5822386176894445
$ws.Range("R7").Value = 5.240147559205
$ws.Range("S7").Value = 0.2258494900888336
$ws.Range("T7").Value = 0.2258494900888337
